# OpenDashboard_Configurations.xlsx - "Adapted option to quit session, bug fixes"
#
# The "BundesagenturFuerArbeit" join-config on TableData is replaced by a
# single, corrected "BundesagenturArbeit" row (the old Nov.2015 filename and
# the separate EW_Altersklassen / Familien rows go away), and the workbook
# now re-opens on TableData (instead of JoinSpec) so the user lands on the
# option to quit/adjust the session right away.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TableData")

# Row 2: correct / rename the Bundesagentur fuer Arbeit entry
$ws.Range("A2").Value = "BundesagenturArbeit"
$ws.Range("B2").Value = "Bundesagentur für Arbeit"
$ws.Range("C2").Value = "BundesagenturArbeit_link.xlsx"
$ws.Range("D2").Value = "BundesagenturArbeit_link"

# Rows 3 & 4 (EW_Altersklassen / Familien) are no longer needed
$ws.Range("A3:D4").ClearContents()

# Make TableData the sheet shown when the workbook is reopened, with A2
# selected (it was JoinSpec before)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
